# "corrected non-admin user in the data"
#
# Rows 2-3 ("To" = Admin) previously also carried a "CC" value (column D)
# wired to the same external user_credentials workbook, but pointed at the
# *Admin* name again via a stray formula -- clear that stray CC value out.
#
# Rows 4-5 ("To" = Admin) were supposed to be the non-admin test rows, but
# their "To" column (C) was wired to the Admin credential; repoint C4/C5 at
# the external user_credentials!$B$4 cell (AutoTestUser) instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the erroneous CC (column D) entries on rows 2 and 3.
$ws.Range("D2").ClearContents()
$ws.Range("D3").ClearContents()

# Rows 4 and 5: point "To" (column C) at the non-admin user instead of admin.
$ws.Range("C4").Formula = "=[1]user_credentials!`$B`$4"
$ws.Range("C5").Formula = "=[1]user_credentials!`$B`$4"
